# Weekly fruit/vegetable price update: a new week's record is inserted
# right after the current row 220, shifting every subsequent row down by
# one (dimension grows from A1:R245 to A1:R246).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 221..245 down to 222..246, duplicating row 221's formatting
# (mirrors Excel's own "Insert Sheet Rows" behaviour) onto the freshly
# opened row 221.
$ws.Rows(221).Insert()

# Populate the newly inserted row 221 with the new week's data.
$ws.Range("A221").Value = 9
$ws.Range("B221").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C221").Value = "Metropolitana"
$ws.Range("D221").Value = 44769
$ws.Range("E221").Value = 13
$ws.Range("F221").Value = 100112026
$ws.Range("G221").Value = "Haba"
$ws.Range("H221").Value = "Sin especificar"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 52
$ws.Range("K221").Value = 20000
$ws.Range("L221").Value = 20000
$ws.Range("M221").Value = 20000
$ws.Range("N221").Value = "$/saco 25 kilos"
$ws.Range("O221").Value = "Región de Coquimbo"
$ws.Range("P221").Value = 800
$ws.Range("Q221").Value = 25
$ws.Range("R221").Value = "Hortaliza"
